$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 1988
$ws.Range("I3").Value = 2116
$ws.Range("H4").Value = 1657
$ws.Range("I4").Value = 530
$ws.Range("I5").Value = 186
$ws.Range("H6").Value = 7922
$ws.Range("I6").Value = 2517
$ws.Range("H7").Value = 25967
$ws.Range("I7").Value = 7337

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 11
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 75
$ws.Range("I3").Value = 70
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 237

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 30
$ws.Range("I6").Value = 39
$ws.Range("I7").Value = 132

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 99
$ws.Range("I4").Value = 17
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 281

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 20
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 46
$ws.Range("I3").Value = 49
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 172

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("H2").Value = 203
$ws.Range("I2").Value = 77
$ws.Range("I7").Value = 253
$ws.Range("I8").Value = 465
$ws.Range("I10").Value = 57
$ws.Range("I11").Value = 126
$ws.Range("I14").Value = 38
$ws.Range("I16").Value = 24
$ws.Range("I19").Value = 211
$ws.Range("I20").Value = 200
$ws.Range("I21").Value = 49
$ws.Range("I22").Value = 21
$ws.Range("I23").Value = 63
$ws.Range("I24").Value = 19
$ws.Range("I28").Value = 3
$ws.Range("I29").Value = 472
$ws.Range("I31").Value = 72
$ws.Range("H33").Value = 1303
$ws.Range("I33").Value = 347
$ws.Range("I37").Value = 237
$ws.Range("I42").Value = 245
$ws.Range("I46").Value = 18
$ws.Range("I50").Value = 28
$ws.Range("I53").Value = 74
$ws.Range("I54").Value = 165
$ws.Range("I55").Value = 82
$ws.Range("I64").Value = 75
$ws.Range("I65").Value = 172
$ws.Range("I67").Value = 281
$ws.Range("I75").Value = 28
$ws.Range("I79").Value = 188
$ws.Range("I83").Value = 138
$ws.Range("I85").Value = 345
$ws.Range("I86").Value = 44
$ws.Range("I88").Value = 58
$ws.Range("H91").Value = 297
$ws.Range("I91").Value = 87
$ws.Range("I93").Value = 38
$ws.Range("I94").Value = 62
$ws.Range("I95").Value = 121
$ws.Range("I99").Value = 132
$ws.Range("H101").Value = 25967
$ws.Range("I101").Value = 7337

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 49
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 123
$ws.Range("H4").Value = 63
$ws.Range("I6").Value = 116
$ws.Range("H7").Value = 1303
$ws.Range("I7").Value = 347

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 39
$ws.Range("I3").Value = 34
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 142
$ws.Range("I3").Value = 156
$ws.Range("I5").Value = 14
$ws.Range("I7").Value = 472

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 85
$ws.Range("I3").Value = 133
$ws.Range("I7").Value = 345

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 64
$ws.Range("I3").Value = 87
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 19

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 18

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 32
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 24
$ws.Range("H7").Value = 297
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 49

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 15
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 61
$ws.Range("I7").Value = 200

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I3").Value = 25
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I3").Value = 29
$ws.Range("H4").Value = 23
$ws.Range("H7").Value = 203
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 142
$ws.Range("I3").Value = 125
$ws.Range("I4").Value = 28
$ws.Range("I6").Value = 155
$ws.Range("I7").Value = 465

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I5").Value = 1
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 253

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I4").Value = 2
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("I2").Value = 1
$ws.Range("I7").Value = 3
